$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# -------------------------------------------------------------------
# Sheet 1 ("Feeling"): recalculated stats for the two existing
# questions (no new shared strings introduced by this part).
# -------------------------------------------------------------------

# Rows 3,5,6 gain an explicit 0.00 number format (value unchanged).
$ws1.Range("C3:D3").NumberFormat = "0.00"
$ws1.Range("C5:D6").NumberFormat = "0.00"

# Rows 4 and 7 keep their existing 0.00 style but get recalculated
# values (and lose the stray yellow highlight fill).
$ws1.Range("C4:D4").ClearFormats()
$ws1.Range("C4:D4").NumberFormat = "0.00"
$ws1.Range("C4").Value = 4.4027777777777777
$ws1.Range("D4").Value = 0.90894230095676964

$ws1.Range("C7:D7").ClearFormats()
$ws1.Range("C7:D7").NumberFormat = "0.00"
$ws1.Range("C7").Value = 4.3305555555555557
$ws1.Range("D7").Value = 0.99534276075614447

# -------------------------------------------------------------------
# Sheet 1, new question rows 11-13: " I find today's camp activities
# difficult." - entered before the Identity sheet content (matches the
# shared-string ordering in the workbook).
# -------------------------------------------------------------------
$ws1.Range("A11").Value = " I find today's camp activities difficult."
$ws1.Range("B11").Value = "Programming Basics"
$ws1.Range("C11").NumberFormat = "0.00"
$ws1.Range("D11").NumberFormat = "0.00"
$ws1.Range("C11").Value = 2.730909090909091
$ws1.Range("D11").Value = 1.2930009127291879

$ws1.Range("A12").Value = " I find today's camp activities difficult."
$ws1.Range("B12").Value = "Micro:bit Pet"
$ws1.Range("C12").NumberFormat = "0.00"
$ws1.Range("D12").NumberFormat = "0.00"
$ws1.Range("C12").Value = 4.0145454545454546
$ws1.Range("D12").Value = 0.91558570803339845

$ws1.Range("A13").Value = " I find today's camp activities difficult."
$ws1.Range("B13").Value = "Technical Design Challenge"
$ws1.Range("C13").NumberFormat = "0.00"
$ws1.Range("D13").NumberFormat = "0.00"
$ws1.Range("C13").Value = 3.0277777777777777
$ws1.Range("D13").Value = 1.4080174550843985

# -------------------------------------------------------------------
# Sheet 2 ("Identity"): entirely new question set with new values.
# -------------------------------------------------------------------

$ws2.Range("A2").Value = "Today's camp activities made me feel like I was a computer scientist."
$ws2.Range("C2").Value = 3.5854545454545454
$ws2.Range("D2").Value = 1.2153846799970776

$ws2.Range("A3").Value = "Today's camp activities made me feel like I was a computer scientist."
$ws2.Range("C3:D3").NumberFormat = "0.00"
$ws2.Range("C3").Value = 3.3927272727272726
$ws2.Range("D3").Value = 1.2722620022659004

$ws2.Range("A4").Value = "Today's camp activities made me feel like I was a computer scientist."
$ws2.Range("C4").Value = 3.7583333333333333
$ws2.Range("D4").Value = 1.3538092343077501

$ws2.Range("A5").Value = "Today's camp activities are useful for what I will be doing in school."
$ws2.Range("C5:D5").NumberFormat = "0.00"
$ws2.Range("C5").Value = 3.6981818181818182
$ws2.Range("D5").Value = 1.130126746507605

$ws2.Range("A6").Value = "Today's camp activities are useful for what I will be doing in school."
$ws2.Range("C6:D6").NumberFormat = "0.00"
$ws2.Range("C6").Value = 3.4763636359999999
$ws2.Range("D6").Value = 1.0796697399999999

$ws2.Range("A7").Value = "Today's camp activities are useful for what I will be doing in school."
$ws2.Range("C7").NumberFormat = "0.00"
$ws2.Range("C7").Font.Color = 0
$ws2.Range("C7").Value = 3.85555556
$ws2.Range("D7").Value = 1.2040887733318328

# New question: "...useful for my future career goals." (rows 8-10)
$ws2.Range("A8").Value = "Today's camp activities are useful for my future career goals."
$ws2.Range("B8").Value = "Programming Basics"
$ws2.Range("C8").NumberFormat = "0.00"
$ws2.Range("D8").NumberFormat = "0.00"
$ws2.Range("C8").Value = 3.8509090909090911
$ws2.Range("D8").Value = 1.1790248239545238

$ws2.Range("A9").Value = "Today's camp activities are useful for my future career goals."
$ws2.Range("B9").Value = "Micro:bit Pet"
$ws2.Range("C9").NumberFormat = "0.00"
$ws2.Range("D9").NumberFormat = "0.00"
$ws2.Range("D9").Font.Color = 0
$ws2.Range("C9").Value = 3.5018181820000001
$ws2.Range("D9").Value = 1.12924611

$ws2.Range("A10").Value = "Today's camp activities are useful for my future career goals."
$ws2.Range("B10").Value = "Technical Design Challenge"
$ws2.Range("C10").NumberFormat = "0.00"
$ws2.Range("D10").NumberFormat = "0.00"
$ws2.Range("C10").Value = 3.8305555555555557
$ws2.Range("D10").Value = 1.2226550962434848

# New question: "I want to do more activities similar..." (rows 11-13)
$ws2.Range("A11").Value = "I want to do more activities similar to today's camp activities."
$ws2.Range("B11").Value = "Programming Basics"
$ws2.Range("C11").NumberFormat = "0.00"
$ws2.Range("D11").NumberFormat = "0.00"
$ws2.Range("C11").Value = 4.2254545454545456
$ws2.Range("D11").Value = 1.0840950650900598

$ws2.Range("A12").Value = "I want to do more activities similar to today's camp activities."
$ws2.Range("B12").Value = "Micro:bit Pet"
$ws2.Range("C12").NumberFormat = "0.00"
$ws2.Range("D12").NumberFormat = "0.00"
$ws2.Range("C12").Value = 3.829090909
$ws2.Range("D12").Value = 1.100864939

$ws2.Range("A13").Value = "I want to do more activities similar to today's camp activities."
$ws2.Range("B13").Value = "Technical Design Challenge"
$ws2.Range("C13").NumberFormat = "0.00"
$ws2.Range("D13").NumberFormat = "0.00"
$ws2.Range("C13").Value = 4.1111111111111107
$ws2.Range("D13").Value = 1.1796269127813648

$ws2.Range("B8:B13").Style = $ws2.Range("B2").Style

# -------------------------------------------------------------------
# Sheet 1, new question rows 8-10: " I find today's camp activities
# interesting. " - entered after the Identity sheet content.
# -------------------------------------------------------------------
$ws1.Range("A8").Value = " I find today's camp activities interesting. "
$ws1.Range("B8").Value = "Programming Basics"
$ws1.Range("C8").NumberFormat = "0.00"
$ws1.Range("D8").NumberFormat = "0.00"
$ws1.Range("C8").Value = 4.3818181818181818
$ws1.Range("D8").Value = 0.98697556694236765

$ws1.Range("A9").Value = " I find today's camp activities interesting. "
$ws1.Range("B9").Value = "Micro:bit Pet"
$ws1.Range("C9").NumberFormat = "0.00"
$ws1.Range("D9").NumberFormat = "0.00"
$ws1.Range("C9").Value = 4.083636363636364
$ws1.Range("D9").Value = 0.8225828652975492

$ws1.Range("A10").Value = " I find today's camp activities interesting. "
$ws1.Range("B10").Value = "Technical Design Challenge"
$ws1.Range("C10").NumberFormat = "0.00"
$ws1.Range("D10").NumberFormat = "0.00"
$ws1.Range("C10").Value = 4.3
$ws1.Range("D10").Value = 1.044417445010805

# -------------------------------------------------------------------
# Sheet 1, new question rows 14-16: " I felt successful after
# completing today's camp activitie" - entered last.
# -------------------------------------------------------------------
$ws1.Range("A14").Value = " I felt successful after completing today's camp activitie"
$ws1.Range("B14").Value = "Programming Basics"
$ws1.Range("C14").NumberFormat = "0.00"
$ws1.Range("D14").NumberFormat = "0.00"
$ws1.Range("C14").Value = 4.3127272727272725
$ws1.Range("D14").Value = 0.99102810021783272

$ws1.Range("A15").Value = " I felt successful after completing today's camp activitie"
$ws1.Range("B15").Value = "Micro:bit Pet"
$ws1.Range("C15").NumberFormat = "0.00"
$ws1.Range("D15").NumberFormat = "0.00"
$ws1.Range("C15").Value = 3.9490909090909092
$ws1.Range("D15").Value = 0.89715301601039432

$ws1.Range("A16").Value = " I felt successful after completing today's camp activitie"
$ws1.Range("B16").Value = "Technical Design Challenge"
$ws1.Range("C16").NumberFormat = "0.00"
$ws1.Range("D16").NumberFormat = "0.00"
$ws1.Range("C16").Value = 4.3111111111111109
$ws1.Range("D16").Value = 1.0303953801480845

$ws1.Range("B8:B16").Style = $ws1.Range("B2").Style

# -------------------------------------------------------------------
# Column widths / selection tweaks observed in the diff.
# -------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 56.83203125

$ws1.Range("E13").Select()
$ws2.Range("C2:D13").Select()
